$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.021.77'
$ws.Range('E2').Value = '  +2.20%  '
$ws.Range('D3').Value = '3.395.17'
$ws.Range('E3').Value = '  +1.82%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '''560.72'
$ws.Range('E5').Value = '  +1.67%  '
$ws.Range('D6').Value = '''174.55'
$ws.Range('E6').Value = '  +2.25%  '
$ws.Range('D7').Value = '''0.625'
$ws.Range('E7').Value = '  +2.73%  '
$ws.Range('D8').Value = '3.384.75'
$ws.Range('E8').Value = '  +1.78%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('E10').Value = '  +11.30%  '
$ws.Range('D11').Value = '''0.631'
$ws.Range('E11').Value = '  +2.56%  '
$ws.Range('D12').Value = '''54.57'
$ws.Range('E12').Value = '  +2.47%  '
$ws.Range('D13').Value = '''0.0000277'
$ws.Range('E13').Value = '  +4.91%  '
$ws.Range('D14').Value = '''9.13'
$ws.Range('E14').Value = '  +2.87%  '
$ws.Range('D15').Value = '3.938.62'
$ws.Range('E15').Value = '  +1.85%  '
$ws.Range('D16').Value = '''18.32'
$ws.Range('E16').Value = '  +3.64%  '
$ws.Range('D17').Value = '3.392.24'
$ws.Range('E17').Value = '  +1.76%  '
$ws.Range('E18').Value = '  +0.73%  '
$ws.Range('D19').Value = '''11.91'
$ws.Range('E19').Value = '  +2.18%  '
$ws.Range('D20').Value = '64.995.30'
$ws.Range('E20').Value = '  +2.33%  '
$ws.Range('D21').Value = '''0.994'
$ws.Range('E21').Value = '  +2.60%  '
$ws.Range('D22').Value = '''473.60'
$ws.Range('E22').Value = '  +17.52%  '
$ws.Range('D23').Value = '''4.98'
$ws.Range('E23').Value = '  +15.76%  '
$ws.Range('D24').Value = '''4.13'
$ws.Range('E24').Value = '  +2.45%  '
$ws.Range('D25').Value = '''86.44'
$ws.Range('E25').Value = '  +4.59%  '
$ws.Range('D26').Value = '''13.71'
$ws.Range('E26').Value = '  +3.49%  '
$ws.Range('D27').Value = '''10.84'
$ws.Range('E27').Value = '  +2.50%  '
$ws.Range('E28').Value = '  +5.02%  '
$ws.Range('D29').Value = '''8.83'
$ws.Range('E29').Value = '  +2.43%  '
$ws.Range('D30').Value = '''30.59'
$ws.Range('E30').Value = '  +4.98%  '
$ws.Range('D31').Value = '''6.69'
$ws.Range('E31').Value = '  +4.80%  '
$ws.Range('D32').Value = '''11.53'
$ws.Range('E32').Value = '  +2.07%  '
$ws.Range('D33').Value = '''585.17'
$ws.Range('E33').Value = '  +1.41%  '
$ws.Range('E34').Value = '  +3.04%  '
$ws.Range('D35').Value = '''60.16'
$ws.Range('E35').Value = '  +4.42%  '
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('B37').Value = 'Stacks'
$ws.Range('C37').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D37').Value = '''3.54'
$ws.Range('E37').Value = '  +3.00%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').Value = '''0.140'
$ws.Range('E38').Value = '  -3.93%  '
$ws.Range('D39').Value = '''35.94'
$ws.Range('E39').Value = '  +1.27%  '
$ws.Range('D40').Value = '0.0₃0749'
$ws.Range('E40').Value = '  +1.66%  '
$ws.Range('D41').Value = '''0.373'
$ws.Range('E41').Value = '  +1.73%  '
$ws.Range('D42').Value = '3.113.35'
$ws.Range('E42').Value = '  -0.80%  '
$ws.Range('D43').Value = '''1.00'
$ws.Range('E43').Value = '  +0.14%  '
$ws.Range('D44').Value = '''2.85'
$ws.Range('E44').Value = '  +1.31%  '
$ws.Range('D45').Value = '''2.51'
$ws.Range('E45').Value = '  +2.64%  '
$ws.Range('D46').Value = '''0.0415'
$ws.Range('E46').Value = '  +3.24%  '
$ws.Range('D47').Value = '''3.21'
$ws.Range('E47').Value = '  +0.95%  '
$ws.Range('E48').Value = '  +4.81%  '
$ws.Range('E49').Value = '  -2.40%  '
$ws.Range('D50').Value = '''136.78'
$ws.Range('E50').Value = '  +3.26%  '
$ws.Range('D51').Value = '''8.32'
$ws.Range('E51').Value = '  +3.51%  '
